# Update "想去人数" (number of people interested) figures that changed
# between data snapshots, on both the "展览" sheet and the combined
# "全部类型" sheet, per the upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5324
$ws1.Range("F3").Value = 378
$ws1.Range("F7").Value = 310
$ws1.Range("F8").Value = 13

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5324
$ws4.Range("F3").Value = 378
$ws4.Range("F8").Value = 310
$ws4.Range("F9").Value = 13

$wb.Save()
